$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for rows 2-25, replacing old Strike# values in column G
$newValues = @{
    2  = 2
    3  = 5
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 3
    9  = 3
    10 = 4
    11 = 6
    12 = 3
    13 = 2
    14 = 4
    15 = 5
    16 = 5
    17 = 5
    18 = 5
    19 = 4
    20 = 2
    21 = 1
    22 = 4
    23 = 4
    24 = 3
    25 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
